$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 149.625
$ws.Cells.Item(9, 9).Value = 148.14285
$ws.Cells.Item(9, 11).Value = 148.14285
$ws.Cells.Item(9, 13).Value = 20.85714999999999
$ws.Cells.Item(17, 8).Value = 2942.2144
$ws.Cells.Item(17, 10).Value = 2942.2144
$ws.Cells.Item(17, 12).Value = 8826.643199999999
$ws.Cells.Item(17, 14).Value = -9162.643199999999
$ws.Cells.Item(43, 8).Value = 3650.6316
$ws.Cells.Item(43, 10).Value = 3427.8
$ws.Cells.Item(43, 12).Value = 3427.8
$ws.Cells.Item(43, 14).Value = -3565.8
$ws.Cells.Item(70, 8).Value = 5395.0527
$ws.Cells.Item(70, 10).Value = 5656.7
$ws.Cells.Item(70, 12).Value = 16970.1
$ws.Cells.Item(70, 14).Value = -17510.1
$ws.Cells.Item(73, 8).Value = 5395.0527
$ws.Cells.Item(73, 10).Value = 5656.7
$ws.Cells.Item(73, 12).Value = 16970.1
$ws.Cells.Item(73, 14).Value = -18842.1
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 13).ClearContents()
$ws.Cells.Item(86, 8).Value = 4292.385
$ws.Cells.Item(86, 9).Value = 4752.6
$ws.Cells.Item(86, 10).Value = 4004.75
$ws.Cells.Item(86, 11).Value = 4752.6
$ws.Cells.Item(86, 12).Value = 4004.75
$ws.Cells.Item(86, 13).Value = -3629.6
$ws.Cells.Item(86, 14).Value = -6250.75
$ws.Cells.Item(89, 8).Value = 4292.385
$ws.Cells.Item(89, 9).Value = 4752.6
$ws.Cells.Item(89, 10).Value = 4004.75
$ws.Cells.Item(89, 11).Value = 23763
$ws.Cells.Item(89, 12).Value = 20023.75
$ws.Cells.Item(89, 13).Value = -18147
$ws.Cells.Item(89, 14).Value = -31255.75
$ws.Cells.Item(109, 8).Value = 44977.355
$ws.Cells.Item(109, 10).Value = 44977.355
$ws.Cells.Item(109, 12).Value = 44977.355
$ws.Cells.Item(109, 14).Value = -47751.355
$ws.Cells.Item(132, 8).Value = 1060.2972
$ws.Cells.Item(132, 9).Value = 1021.82355
$ws.Cells.Item(132, 10).Value = 1496.3334
$ws.Cells.Item(132, 11).Value = 3065.47065
$ws.Cells.Item(132, 12).Value = 4489.0002
$ws.Cells.Item(132, 13).Value = -535.4706499999998
$ws.Cells.Item(132, 14).Value = -9549.0002
$ws.Cells.Item(137, 8).Value = 3729.0952
$ws.Cells.Item(137, 10).Value = 7532
$ws.Cells.Item(137, 12).Value = 22596
$ws.Cells.Item(137, 14).Value = -27696
$ws.Cells.Item(138, 8).Value = 5181.7837
$ws.Cells.Item(138, 9).Value = 3594.1428
$ws.Cells.Item(138, 10).Value = 5552.2334
$ws.Cells.Item(138, 11).Value = 10782.4284
$ws.Cells.Item(138, 12).Value = 16656.7002
$ws.Cells.Item(138, 13).Value = -5642.428400000001
$ws.Cells.Item(138, 14).Value = -26936.7002
$ws.Cells.Item(141, 8).Value = 6139.5
$ws.Cells.Item(141, 10).Value = 4999.5
$ws.Cells.Item(141, 12).Value = 14998.5
$ws.Cells.Item(141, 14).Value = -25358.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16384.334
$ws.Cells.Item(32, 9).Value = 12496.677
$ws.Cells.Item(32, 10).Value = 32906.875
$ws.Cells.Item(32, 11).Value = 12496.677
$ws.Cells.Item(32, 12).Value = 32906.875
$ws.Cells.Item(32, 13).Value = -12209.677
$ws.Cells.Item(32, 14).Value = -33480.875
$ws.Cells.Item(74, 8).Value = 2377.3076
$ws.Cells.Item(74, 9).Value = 2250.5
$ws.Cells.Item(74, 10).Value = 3074.75
$ws.Cells.Item(74, 11).Value = 2250.5
$ws.Cells.Item(74, 12).Value = 3074.75
$ws.Cells.Item(74, 13).Value = -1376.5
$ws.Cells.Item(74, 14).Value = -4822.75
$ws.Cells.Item(77, 8).Value = 2377.3076
$ws.Cells.Item(77, 9).Value = 2250.5
$ws.Cells.Item(77, 10).Value = 3074.75
$ws.Cells.Item(77, 11).Value = 11252.5
$ws.Cells.Item(77, 12).Value = 15373.75
$ws.Cells.Item(77, 13).Value = -6884.5
$ws.Cells.Item(77, 14).Value = -24109.75
$ws.Cells.Item(110, 8).Value = 4354
$ws.Cells.Item(110, 9).Value = 4354
$ws.Cells.Item(110, 11).Value = 4354
$ws.Cells.Item(110, 13).Value = -2309
$ws.Cells.Item(132, 8).Value = 5161.5674
$ws.Cells.Item(132, 9).Value = 4820.643
$ws.Cells.Item(132, 11).Value = 14461.929
$ws.Cells.Item(132, 13).Value = -11931.929
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 42593.332
$ws.Cells.Item(81, 10).Value = 42593.332
$ws.Cells.Item(81, 12).Value = 42593.332
$ws.Cells.Item(81, 14).Value = -44715.332
$ws.Cells.Item(84, 8).Value = 42593.332
$ws.Cells.Item(84, 10).Value = 42593.332
$ws.Cells.Item(84, 12).Value = 127779.996
$ws.Cells.Item(84, 14).Value = -138387.996
$ws.Cells.Item(94, 8).Value = 6062549.5
$ws.Cells.Item(94, 9).Value = 1561.75
$ws.Cells.Item(94, 11).Value = 1561.75
$ws.Cells.Item(94, 13).Value = -1110.75
$ws.Cells.Item(99, 8).Value = 1312.9166
$ws.Cells.Item(99, 9).Value = 1126
$ws.Cells.Item(99, 11).Value = 1126
$ws.Cells.Item(99, 13).Value = 372
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1095.0834
$ws.Cells.Item(22, 9).Value = 1077.2858
$ws.Cells.Item(22, 11).Value = 1077.2858
$ws.Cells.Item(22, 13).Value = -727.2858000000001
$ws.Cells.Item(31, 8).Value = 3582.319
$ws.Cells.Item(31, 9).Value = 2606.3
$ws.Cells.Item(31, 10).Value = 5304.706
$ws.Cells.Item(31, 11).Value = 2606.3
$ws.Cells.Item(31, 12).Value = 5304.706
$ws.Cells.Item(31, 13).Value = -2311.3
$ws.Cells.Item(31, 14).Value = -5894.706
$ws.Cells.Item(34, 8).Value = 3582.319
$ws.Cells.Item(34, 9).Value = 2606.3
$ws.Cells.Item(34, 10).Value = 5304.706
$ws.Cells.Item(34, 11).Value = 2606.3
$ws.Cells.Item(34, 12).Value = 5304.706
$ws.Cells.Item(34, 13).Value = -2404.3
$ws.Cells.Item(34, 14).Value = -5708.706
$ws.Cells.Item(132, 8).Value = 5461.4287
$ws.Cells.Item(132, 9).Value = 4921.3335
$ws.Cells.Item(132, 11).Value = 14764.0005
$ws.Cells.Item(132, 13).Value = -12234.0005
$ws.Cells.Item(141, 8).Value = 342632.72
$ws.Cells.Item(141, 10).Value = 342632.72
$ws.Cells.Item(141, 12).Value = 342632.72
$ws.Cells.Item(141, 14).Value = -352992.72
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 783.6
$ws.Cells.Item(11, 9).Value = 690.8570999999999
$ws.Cells.Item(11, 11).Value = 2072.5713
$ws.Cells.Item(11, 13).Value = -1932.5713
$ws.Cells.Item(12, 8).Value = 193.43478
$ws.Cells.Item(12, 9).Value = 196.25
$ws.Cells.Item(12, 10).Value = 192.8421
$ws.Cells.Item(12, 11).Value = 588.75
$ws.Cells.Item(12, 12).Value = 578.5263
$ws.Cells.Item(12, 13).Value = -415.75
$ws.Cells.Item(12, 14).Value = -924.5263
$ws.Cells.Item(68, 8).Value = 1227
$ws.Cells.Item(68, 10).Value = 1227
$ws.Cells.Item(68, 12).Value = 3681
$ws.Cells.Item(68, 14).Value = -5303
$ws.Cells.Item(71, 8).Value = 1227
$ws.Cells.Item(71, 10).Value = 1227
$ws.Cells.Item(71, 12).Value = 11043
$ws.Cells.Item(71, 14).Value = -19155
$ws.Cells.Item(113, 8).Value = 2259.2942
$ws.Cells.Item(113, 10).Value = 2092.2
$ws.Cells.Item(113, 12).Value = 6276.599999999999
$ws.Cells.Item(113, 14).Value = -10616.6
$ws.Cells.Item(131, 8).Value = 3847.5518
$ws.Cells.Item(131, 10).Value = 4595.4
$ws.Cells.Item(131, 12).Value = 13786.2
$ws.Cells.Item(131, 14).Value = -23866.2
$ws.Cells.Item(132, 8).Value = 4014
$ws.Cells.Item(132, 9).Value = 1524.75
$ws.Cells.Item(132, 11).Value = 13722.75
$ws.Cells.Item(132, 13).Value = -11192.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 41808736
$ws.Cells.Item(80, 10).Value = 166668670
$ws.Cells.Item(80, 12).Value = 166668670
$ws.Cells.Item(80, 14).Value = -166670666
$ws.Cells.Item(83, 8).Value = 41808736
$ws.Cells.Item(83, 10).Value = 166668670
$ws.Cells.Item(83, 12).Value = 833343350
$ws.Cells.Item(83, 14).Value = -833353334
$ws.Cells.Item(122, 8).Value = 5208.7144
$ws.Cells.Item(122, 9).Value = 2947
$ws.Cells.Item(122, 10).Value = 6905
$ws.Cells.Item(122, 11).Value = 8841
$ws.Cells.Item(122, 12).Value = 20715
$ws.Cells.Item(122, 13).Value = -6391
$ws.Cells.Item(122, 14).Value = -25615
$ws.Cells.Item(132, 8).Value = 3745.5688
$ws.Cells.Item(132, 9).Value = 3932.8936
$ws.Cells.Item(132, 10).Value = 2945.182
$ws.Cells.Item(132, 11).Value = 11798.6808
$ws.Cells.Item(132, 12).Value = 8835.545999999998
$ws.Cells.Item(132, 13).Value = -9268.6808
$ws.Cells.Item(132, 14).Value = -13895.546
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 14040.077
$ws.Cells.Item(40, 9).Value = 21051.834
$ws.Cells.Item(40, 11).Value = 21051.834
$ws.Cells.Item(40, 13).Value = -20915.834
$ws.Cells.Item(93, 8).Value = 13002513
$ws.Cells.Item(93, 9).Value = 2926.7693
$ws.Cells.Item(93, 10).Value = 37144600
$ws.Cells.Item(93, 11).Value = 2926.7693
$ws.Cells.Item(93, 12).Value = 37144600
$ws.Cells.Item(93, 13).Value = -1678.7693
$ws.Cells.Item(93, 14).Value = -37147096
$ws.Cells.Item(100, 8).Value = 72527.31
$ws.Cells.Item(100, 9).Value = 160086.14
$ws.Cells.Item(100, 11).Value = 160086.14
$ws.Cells.Item(100, 13).Value = -159545.14
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 5156.3125
$ws.Cells.Item(122, 9).Value = 4323.727
$ws.Cells.Item(122, 10).Value = 6988
$ws.Cells.Item(122, 11).Value = 12971.181
$ws.Cells.Item(122, 12).Value = 20964
$ws.Cells.Item(122, 13).Value = -10521.181
$ws.Cells.Item(122, 14).Value = -25864
$ws.Cells.Item(132, 8).Value = 3117.8572
$ws.Cells.Item(132, 9).Value = 1868
$ws.Cells.Item(132, 11).Value = 5604
$ws.Cells.Item(132, 13).Value = -3074
$ws.Cells.Item(136, 8).Value = 3161.348
$ws.Cells.Item(136, 9).Value = 2422.3125
$ws.Cells.Item(136, 10).Value = 4850.5713
$ws.Cells.Item(136, 11).Value = 7266.9375
$ws.Cells.Item(136, 12).Value = 14551.7139
$ws.Cells.Item(136, 13).Value = -4716.9375
$ws.Cells.Item(136, 14).Value = -19651.7139
